$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B, C, D, F, G, I, K, L across rows 2-25
# Row format: row number, then values for B, C, D, F, G, I, K, L in order
$data = @(
    ,(2, 14.41155629258586, 4.734212120649882, 8.311726735119583, 42.41018029824375, 3.724044295903637, 34.83494515647468, 12.45006114069616, 10.89282799992328)
    ,(3, 14.29492639753532, 4.548862608271972, 8.303148054274754, 42.05091988384387, 3.727189311259263, 34.65730748989031, 12.37306244056855, 10.87242367938386)
    ,(4, 14.22793586449322, 4.429973538784413, 8.297682455629513, 41.8354336038399, 3.729220532384572, 34.55198532402104, 12.32963850521564, 10.86216072879885)
    ,(5, 14.20183014944285, 4.380279702728996, 8.295405010251887, 41.74896793574553, 3.730073550832025, 34.51002852366478, 12.31292986660309, 10.85855070480773)
    ,(6, 14.19756827111841, 4.371953896831874, 8.295023802704023, 41.73469353383176, 3.730216723380179, 34.50312034955974, 12.31021550979405, 10.85798588882716)
    ,(7, 14.22757892156468, 4.429308346947451, 8.297651944774371, 41.83426196231946, 3.729231934014545, 34.55141555570083, 12.32940914811589, 10.86210972285725)
    ,(8, 14.37040296916304, 4.671377007809244, 8.308809762906726, 42.28528780203904, 3.725107962036151, 34.77293037673358, 12.42272349152347, 10.88532408520837)
    ,(9, 14.68563950205059, 5.104518062768327, 8.329131592760412, 43.20729412418424, 3.717811506482786, 35.23611188264815, 12.63541599322527, 10.94868993567033)
    ,(10, 14.93651619518809, 5.396060653680233, 8.343132043743427, 43.90352067126772, 3.712926978826145, 35.5926217903194, 12.80850878383533, 11.00591659278993)
    ,(11, 15.05433704129481, 5.522678768854877, 8.349303545533196, 44.22344357175569, 3.710807037071446, 35.75801834239349, 12.89060576403404, 11.03421434520382)
    ,(12, 15.09944108752209, 5.569747585420084, 8.351612474444249, 44.34497393091723, 3.710018850664735, 35.82108464831722, 12.92214942229642, 11.04525040608166)
    ,(13, 15.08970612857927, 5.559649767968184, 8.351116449676763, 44.31878443621442, 3.710187953197504, 35.80748334295795, 12.91533612482782, 11.04285943732496)
    ,(14, 15.05803825336247, 5.526568836837587, 8.349494063678568, 44.23343473288094, 3.710741900595675, 35.76319830965705, 12.89319191410638, 11.0351159014213)
    ,(15, 15.03870299138788, 5.506190972683481, 8.348496656684876, 44.18120305034958, 3.711083106821135, 35.73612813350434, 12.87968641557739, 11.030414309093)
    ,(16, 14.92888698488879, 5.387663660044124, 8.342724775415315, 43.8826707691679, 3.713067568258599, 35.58187507603449, 12.80320877509289, 11.00411239563098)
    ,(17, 14.86243423495536, 5.313401169093813, 8.339133573728031, 43.7002963165808, 3.714311048473235, 35.488050079877, 12.75713269097367, 10.98855346029734)
    ,(18, 14.82456290482007, 5.270122755419824, 8.337049414228284, 43.59570585794211, 3.71503587634586, 35.43438912530974, 12.73094854585197, 10.97981792181985)
    ,(19, 14.81180184319535, 5.255372959264781, 8.336340550758464, 43.56034835153957, 3.715282943917603, 35.41627363722616, 12.72213840917618, 10.97689705312808)
    ,(20, 14.86947227820377, 5.321365088204782, 8.33951778566764, 43.71967923394941, 3.714177683811651, 35.49800657987118, 12.76200489342281, 10.99018767166656)
    ,(21, 15.06732698685929, 5.536309463881336, 8.349971357637562, 44.25849425648781, 3.710578797542221, 35.77619434479964, 12.89968407040655, 11.03738171750314)
    ,(22, 15.19945937713807, 5.671660929398348, 8.356639806391259, 44.61283389444155, 3.708311716771646, 35.96052326715314, 12.99230604952208, 11.07008999559013)
    ,(23, 15.128694110039, 5.599894776674679, 8.35309559875423, 44.42354154029059, 3.709513951212134, 35.86192279214939, 12.94263956982259, 11.05246428838893)
    ,(24, 14.86628933754024, 5.317766418027579, 8.339344144295778, 43.71091540407509, 3.714237947022655, 35.49350437094611, 12.75980121669443, 10.98944819174697)
    ,(25, 14.59681175599301, 4.991944160585953, 8.323798760493466, 42.95425947538171, 3.719701350595489, 35.10786042630814, 12.57482352951874, 10.85855070480773)
)

$cols = @("B", "C", "D", "F", "G", "I", "K", "L")

foreach ($entry in $data) {
    $row = $entry[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $col = $cols[$i]
        $value = $entry[$i + 1]
        $ws.Range("$col$row").Value = $value
    }
}

